$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns hold plain text (not real numbers) in the source data,
# e.g. "25.898.03", "1.890", "  -0.26%  ". Force each target cell to Text format
# before assigning so Excel does not reinterpret numeric-looking values as
# numbers (which would also silently drop trailing zeros such as in "1.890").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.898.03"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.631.77"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "215.91"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "0.5105"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.2569"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "0.06336"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "0.07781"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "4.269"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.640.27"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "1.857.67"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "0.5498"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "63.78"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "0.0₅7632"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "25.930.12"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "4.413"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "194.38"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "9.849"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "6.025"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "1.890"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "141.97"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "0.1251"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("D28").Value = "6.748"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "0.04883"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "3.232"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "3.181"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "0.8974"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "2.541"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "0.5497"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "1.116.20"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").Value = "0.01558"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "5.577"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "0.7961"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "97.54"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "1.766.67"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  -8.72%  "
$ws.Range("D47").Value = "0.4440"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "54.67"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "0.05130"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "7.546"
$ws.Range("E51").Value = "  +2.64%  "
